# Insercao do link do GitHub.
# Slide 1, Shape 58 ("Shape 58") - last paragraph currently reads
# "*falta colocar o GitHub*". Replace it with the project's GitHub URL as a
# hyperlink run, followed by a plain space run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$url = "https://github.com/FernandoLins8/Projeto-Estruturas-de-dados"

$para = $tr.Paragraphs(5, 1)

# Replace the paragraph text with the URL plus a trailing space - this keeps
# the paragraph's own pPr/endParaRPr intact while giving us a single run to
# split below.
$para.Text = $url + " "

# Hyperlink only the URL portion (not the trailing space) - addressing it as
# its own character range keeps it as a separate run from the space run.
$urlRange = $para.Characters(1, $url.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $url

Write-Output "Inserted GitHub hyperlink into slide 1."
